# Updated files and excel
#
# 1. LoginCredentials (sheet 1): the previous selection/active-cell marker
#    (G5) is cleared and the cursor is left on A2 instead; the sheet is no
#    longer the active tab once GroupDetails is added.
# 2. A new worksheet, "GroupDetails", is added right after LoginCredentials
#    and becomes the active sheet/tab.
# 3. GroupDetails gets a small 2x2 table of data (Location/Groups headers
#    with a Dallas/Test123 row) and the cursor is left on G14.

$wb = $excel.ActiveWorkbook

# --- LoginCredentials: move the selection from G5 to A2 --------------------
$loginSheet = $wb.Worksheets.Item(1)
$loginSheet.Activate()
$loginSheet.Range("A2").Select()

# --- Add the GroupDetails worksheet right after LoginCredentials -----------
$groupSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $loginSheet)
$groupSheet.Name = "GroupDetails"

# --- Populate GroupDetails --------------------------------------------------
$groupSheet.Range("A1").Value = "Location"
$groupSheet.Range("B1").Value = "Groups"
$groupSheet.Range("A2").Value = "Dallas"
$groupSheet.Range("B2").Value = "Test123"

# Leave the cursor on G14, matching the saved selection in the workbook.
$groupSheet.Range("G14").Select()
